$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Part 1: refresh time_taken timestamps on the "data" sheet ---
$timestamps = @(
    "2021-10-05 14:20:13.368320",
    "2021-10-05 14:20:13.368328",
    "2021-10-05 14:20:13.368331",
    "2021-10-05 14:20:13.368334",
    "2021-10-05 14:20:13.368336",
    "2021-10-05 14:20:13.368339",
    "2021-10-05 14:20:13.368341",
    "2021-10-05 14:20:13.368344",
    "2021-10-05 14:20:13.368347",
    "2021-10-05 14:20:13.368349",
    "2021-10-05 14:20:13.368351",
    "2021-10-05 14:20:13.368354",
    "2021-10-05 14:20:13.368357",
    "2021-10-05 14:20:13.368359",
    "2021-10-05 14:20:13.368362",
    "2021-10-05 14:20:13.368364",
    "2021-10-05 14:20:13.368367",
    "2021-10-05 14:20:13.368370",
    "2021-10-05 14:20:13.368372",
    "2021-10-05 14:20:13.368375",
    "2021-10-05 14:20:13.368378",
    "2021-10-05 14:20:13.368380",
    "2021-10-05 14:20:13.368383"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- Part 2: add a new "metadata" sheet after "data" ---
$new = $wb.Worksheets.Add($null, $dataSheet)
$new.Name = "metadata"

$new.Range("B1").Value = "data_name"
$new.Range("C1").Value = "data_id"
$new.Range("D1").Value = "data_version"
$new.Range("E1").Value = "data_version_created"
$new.Range("F1").Value = "panel_query_time"
$new.Range("G1").Value = "panel_get_request"

$new.Range("A2").Value = 0
$new.Range("B2").Value = "Familial dysautonomia"
$new.Range("C2").Value = 7
$new.Range("D2").NumberFormat = "@"
$new.Range("D2").Value = "1.15"
$new.Range("D2").ClearFormats()
$new.Range("E2").Value = "2021-09-14T10:32:22.094881Z"
$new.Range("F2").Value = "2021-10-05 14:20:13.364711"
$new.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/7/?format=json"

# Copy header + index-column formatting from the "data" sheet so the new
# sheet's style matches (bold, centered, bordered header style).
$dataSheet.Range("B1:F1").Copy()
$new.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$new.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$dataSheet.Activate()

Write-Output "done"
